$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet by copying the existing "2022-Q2" sheet
#    (so it inherits the identical column layout / styles), place it right
#    after "总计", rename it, and fill in the new quarter's figures.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")
$wsQ2 = $wb.Worksheets.Item("2022-Q2")

$wsQ2.Copy($null, $wsTotal)
$wsQ3 = $wb.Worksheets.Item(2)
$wsQ3.Name = "2022-Q3"

# Row 2 holds the fund's figures for this quarter. Keep the fund identity
# (A2, B2, C2) as-is (copied from 2022-Q2) and only update the metrics.
$wsQ3.Range("D2").Value = "'12.19"
$wsQ3.Range("D2").Style = "Normal"

$wsQ3.Range("E2").Value = "'99.43"
$wsQ3.Range("E2").Style = "Normal"

$wsQ3.Range("F2").Value = "'2.56"
$wsQ3.Range("F2").Style = "Normal"

$wsQ3.Range("G2").Value = "'0.3121"
$wsQ3.Range("G2").Style = "Normal"

$wsQ3.Range("H2").Value = 10

# ---------------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: insert a new row right below the
#    header for 2022-Q3 and push the existing quarterly rows down by one.
# ---------------------------------------------------------------------------
$wsTotal.Rows("2:2").Insert()

# B2:D2 should use the plain (unstyled) formatting used by every other data
# row, not the bold header formatting Excel copied down from row 1.
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("B2").Style = "Normal"
$wsTotal.Range("C2").Value = 1
$wsTotal.Range("C2").Style = "Normal"
$wsTotal.Range("D2").Value = 0.31
$wsTotal.Range("D2").Style = "Normal"

# Column A is the running index and uses the same (bold/bordered) style as
# the rest of the index column; copy that formatting down from row 3 (the
# old row 2, which still has it) onto the freshly inserted row 2.
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)

# Renumber the running index in column A (0,1,2,...) for every data row now
# that an additional row has been inserted.
$wsTotal.Range("A2").Value = 0
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("A4").Value = 2
$wsTotal.Range("A5").Value = 3
$wsTotal.Range("A6").Value = 4
$wsTotal.Range("A7").Value = 5
$wsTotal.Range("A8").Value = 6
$wsTotal.Range("A9").Value = 7
